# Append the 2025-03-14 price row (row 13) to every price-history sheet,
# carrying forward the same value that was recorded on 2025-03-13 (row 12).

$wb = $excel.ActiveWorkbook

$newDate = "2025-03-14"

$sheetValues = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "43"
    "N-type Wafer"               = "1.19"
    "Cell Topcon 183mm"          = "0.298"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,335"
    "Silver Busbar front-side"   = "7,987"
    "Silver finger front-side"   = "8,037"
    "USD_CNY"                    = "7.2567"
}

foreach ($sheetName in $sheetValues.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $value = $sheetValues[$sheetName]

    $dateCell = $ws.Range("A13")
    $priceCell = $ws.Range("B13")

    # Force text storage (matching the rest of the column) instead of
    # letting the date-looking / numeric-looking strings get auto-typed,
    # then clear the format back to the sheet's default style so no new
    # cell style gets introduced.
    $dateCell.NumberFormat = "@"
    $priceCell.NumberFormat = "@"

    $dateCell.Value = $newDate
    $priceCell.Value = $value

    $ws.Range("A13:B13").ClearFormats()
}
